$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (price + volume/1h) per the commit diff.
# Values that look like plain decimal numbers are prefixed with a leading
# apostrophe so Excel stores them as text (matching the source t="inlineStr"
# cells) instead of silently coercing them to numeric cells.

$ws.Range("D2").Value = "72.949.54"
$ws.Range("E2").Value = "  +1.15%  "
$ws.Range("D3").Value = "4.026.00"
$ws.Range("E3").Value = "  +0.12%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'594.59"
$ws.Range("E5").Value = "  +12.34%  "
$ws.Range("D6").Value = "'153.25"
$ws.Range("E6").Value = "  +1.28%  "
$ws.Range("D7").Value = "'0.688"
$ws.Range("E7").Value = "  -2.52%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "'0.762"
$ws.Range("E9").Value = "  +1.54%  "
$ws.Range("E10").Value = "  -1.13%  "
$ws.Range("D11").Value = "'54.32"
$ws.Range("E11").Value = "  +11.88%  "
$ws.Range("D12").Value = "'0.0000321"
$ws.Range("E12").Value = "  -1.84%  "
$ws.Range("D13").Value = "'10.99"
$ws.Range("E13").Value = "  +3.23%  "
$ws.Range("D14").Value = "4.674.20"
$ws.Range("E14").Value = "  +0.23%  "
$ws.Range("D15").Value = "4.030.17"
$ws.Range("E15").Value = "  -0.14%  "
$ws.Range("D16").Value = "'1.27"
$ws.Range("E16").Value = "  +6.49%  "
$ws.Range("D17").Value = "'14.30"
$ws.Range("E17").Value = "  +1.55%  "
$ws.Range("D18").Value = "'20.68"
$ws.Range("E18").Value = "  +0.71%  "
$ws.Range("E19").Value = "  -0.40%  "
$ws.Range("D20").Value = "72.827.69"
$ws.Range("E20").Value = "  +1.10%  "
$ws.Range("D21").Value = "'439.86"
$ws.Range("E21").Value = "  +2.53%  "
$ws.Range("D22").Value = "'4.75"
$ws.Range("E22").Value = "  +12.74%  "
$ws.Range("D23").Value = "'97.38"
$ws.Range("E23").Value = "  -0.91%  "
$ws.Range("E24").Value = "  +1.97%  "
$ws.Range("D25").Value = "'14.35"
$ws.Range("E25").Value = "  +0.71%  "
$ws.Range("D26").Value = "'4.32"
$ws.Range("E26").Value = "  +19.56%  "
$ws.Range("D27").Value = "'11.55"
$ws.Range("E27").Value = "  +1.44%  "
$ws.Range("D28").Value = "'10.76"
$ws.Range("E28").Value = "  +0.41%  "
$ws.Range("D29").Value = "'5.93"
$ws.Range("E29").Value = "  +1.43%  "
$ws.Range("D30").Value = "'36.79"
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("D31").Value = "'7.95"
$ws.Range("E31").Value = "  +10.37%  "
$ws.Range("B32").Value = "Cosmos"
$ws.Range("C32").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D32").Value = "'13.66"
$ws.Range("E32").Value = "  +1.97%  "
$ws.Range("B33").Value = "InjectiveProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D33").Value = "'50.28"
$ws.Range("E33").Value = "  +11.35%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.133"
$ws.Range("E34").Value = "  +2.28%  "
$ws.Range("D35").Value = "'690.97"
$ws.Range("E35").Value = "  +1.13%  "
$ws.Range("D36").Value = "'71.40"
$ws.Range("E36").Value = "  +8.71%  "
$ws.Range("D37").Value = "'0.447"
$ws.Range("E37").Value = "  -0.32%  "
$ws.Range("D38").Value = "0.0₃0874"
$ws.Range("E38").Value = "  +5.48%  "
$ws.Range("D39").Value = "'3.44"
$ws.Range("E39").Value = "  +5.21%  "
$ws.Range("E40").Value = "  -0.73%  "
$ws.Range("D41").Value = "'11.21"
$ws.Range("E41").Value = "  +12.98%  "
$ws.Range("D42").Value = "'3.36"
$ws.Range("E42").Value = "  -1.86%  "
$ws.Range("D43").Value = "'0.999"
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("E44").Value = "  +1.24%  "
$ws.Range("D45").Value = "'1.00"
$ws.Range("E45").Value = "  +0.14%  "
$ws.Range("E46").Value = "  +0.33%  "
$ws.Range("D47").Value = "'2.75"
$ws.Range("E47").Value = "  +1.14%  "
$ws.Range("D48").Value = "'3.38"
$ws.Range("E48").Value = "  -0.93%  "
$ws.Range("D49").Value = "'3.51"
$ws.Range("E49").Value = "  +7.99%  "
$ws.Range("D50").Value = "'3.05"
$ws.Range("E50").Value = "  +0.95%  "
$ws.Range("D51").Value = "'2.17"
$ws.Range("E51").Value = "  +8.46%  "
